$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Contact / delete) return-example column: replace the generic
# placeholder return example with the concrete one for deleting a contact /
# rejecting a friend request.
$ws.Range("D11").Value = "[{Status=Success},{username=user_name},{contactusername=contact}]"

# Move the active selection from E16 to D16, matching the saved view state.
$ws.Range("D16").Select()
